$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

# ---------------------------------------------------------------------------
# Cell values (rows 3-9 are new contributor rows; row 2 keeps its data but
# gains new formatting below). Columns: A = Sr No, B = Name, C = Roles,
# D = Fields, E = Hours per week.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Prajwala Pandit"
$ws.Range("C3").Value = "Controller"
$ws.Range("D3").Value = "Requirements" + $nl + "Ana and design" + $nl + "Dev and V"
$ws.Range("E3").Value = 5

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Payal Shah"
$ws.Range("C4").Value = "Controller"
$ws.Range("D4").Value = "Ana and design" + $nl + "Development and Ver"
$ws.Range("E4").Value = 5

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Harsimran singh"
$ws.Range("C5").Value = "controller and DSP"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Bharat Bansal"
$ws.Range("C6").Value = "Controller & client app"
$ws.Range("D6").Value = "Analysis and design , Dev"
$ws.Range("E6").Value = 5

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Danish Ahmed"
$ws.Range("C7").Value = "Controller"
$ws.Range("D7").Value = "Ana and design" + $nl + "Development and Ver"
$ws.Range("E7").Value = 5

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Raghuraj"
$ws.Range("C8").Value = "Controller"
$ws.Range("D8").Value = "Analysis and Design, Development and verification"
$ws.Range("E8").Value = 5

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Adit"
$ws.Range("C9").Value = "Controller & client app"
$ws.Range("D9").Value = "Requirements" + $nl + "Ana and design" + $nl + "Dev and V"
$ws.Range("E9").Value = 5

# ---------------------------------------------------------------------------
# Column widths (B & C widen a bit to fit the new longer entries)
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 16.28515625
$ws.Columns("C").ColumnWidth = 22.42578125

# ---------------------------------------------------------------------------
# Header row now gets a bottom divider (thick bottom border look via row)
# ---------------------------------------------------------------------------
$ws.Rows(1).RowHeight = 19.5

# ---------------------------------------------------------------------------
# Row 2 (existing entry) - add the thin grey "card" border + right-aligned
# numeric columns + wrap-texted role/fields columns.
# ---------------------------------------------------------------------------
$ws.Range("A2").Font.Size = 12.1
$ws.Range("A2").Font.Color = 0
$ws.Range("A2").Borders.LineStyle = 1
$ws.Range("A2").Borders.Weight = -4138
$ws.Range("A2").Borders.Color = 13421772
$ws.Range("A2").HorizontalAlignment = -4152

$ws.Range("B2").Font.Size = 12.1
$ws.Range("B2").Font.Color = 0
$ws.Range("B2").Borders.LineStyle = 1
$ws.Range("B2").Borders.Weight = -4138
$ws.Range("B2").Borders.Color = 13421772

$ws.Range("C2").Font.Size = 12.1
$ws.Range("C2").Font.Color = 0
$ws.Range("C2").Borders.LineStyle = 1
$ws.Range("C2").Borders.Weight = -4138
$ws.Range("C2").Borders.Color = 13421772
$ws.Range("C2").WrapText = $true

$ws.Range("D2").Font.Size = 12.1
$ws.Range("D2").Font.Color = 0
$ws.Range("D2").Borders.LineStyle = 1
$ws.Range("D2").Borders.Weight = -4138
$ws.Range("D2").Borders.Color = 13421772
$ws.Range("D2").WrapText = $true

$ws.Range("E2").Font.Size = 12.1
$ws.Range("E2").Font.Color = 0
$ws.Range("E2").Borders.LineStyle = 1
$ws.Range("E2").Borders.Weight = -4138
$ws.Range("E2").Borders.Color = 13421772
$ws.Range("E2").HorizontalAlignment = -4152

$ws.Rows(2).RowHeight = 63.75

# ---------------------------------------------------------------------------
# Row 3
# ---------------------------------------------------------------------------
$ws.Range("A3").Font.Size = 12.1
$ws.Range("A3").Borders.LineStyle = 1
$ws.Range("A3").Borders.Weight = -4138
$ws.Range("A3").Borders.Color = 13421772
$ws.Range("A3").HorizontalAlignment = -4152

$ws.Range("B3").Font.Size = 12.1
$ws.Range("B3").Borders.LineStyle = 1
$ws.Range("B3").Borders.Weight = -4138
$ws.Range("B3").Borders.Color = 13421772

$ws.Range("C3").Font.Size = 12.1
$ws.Range("C3").Borders.LineStyle = 1
$ws.Range("C3").Borders.Weight = -4138
$ws.Range("C3").Borders.Color = 13421772

$ws.Range("D3").Font.Size = 12.1
$ws.Range("D3").Font.Color = 0
$ws.Range("D3").Borders.LineStyle = 1
$ws.Range("D3").Borders.Weight = -4138
$ws.Range("D3").Borders.Color = 13421772
$ws.Range("D3").WrapText = $true

$ws.Range("E3").Font.Size = 12.1
$ws.Range("E3").Borders.LineStyle = 1
$ws.Range("E3").Borders.Weight = -4138
$ws.Range("E3").Borders.Color = 13421772
$ws.Range("E3").HorizontalAlignment = -4152

$ws.Rows(3).RowHeight = 48

# ---------------------------------------------------------------------------
# Row 4
# ---------------------------------------------------------------------------
$ws.Range("A4").Font.Size = 12.1
$ws.Range("A4").Borders.LineStyle = 1
$ws.Range("A4").Borders.Weight = -4138
$ws.Range("A4").Borders.Color = 13421772
$ws.Range("A4").HorizontalAlignment = -4152

$ws.Range("B4").Font.Size = 12.1
$ws.Range("B4").Borders.LineStyle = 1
$ws.Range("B4").Borders.Weight = -4138
$ws.Range("B4").Borders.Color = 13421772

$ws.Range("C4").Font.Size = 12.1
$ws.Range("C4").Borders.LineStyle = 1
$ws.Range("C4").Borders.Weight = -4138
$ws.Range("C4").Borders.Color = 13421772

$ws.Range("D4").Font.Size = 12.1
$ws.Range("D4").Borders.LineStyle = 1
$ws.Range("D4").Borders.Weight = -4138
$ws.Range("D4").Borders.Color = 13421772
$ws.Range("D4").WrapText = $true

$ws.Range("E4").Font.Size = 12.1
$ws.Range("E4").Borders.LineStyle = 1
$ws.Range("E4").Borders.Weight = -4138
$ws.Range("E4").Borders.Color = 13421772
$ws.Range("E4").HorizontalAlignment = -4152

$ws.Rows(4).RowHeight = 48

# ---------------------------------------------------------------------------
# Row 5 (no fields/hours filled in for Harsimran)
# ---------------------------------------------------------------------------
$ws.Range("A5").Font.Size = 12.1
$ws.Range("A5").Borders.LineStyle = 1
$ws.Range("A5").Borders.Weight = -4138
$ws.Range("A5").Borders.Color = 13421772
$ws.Range("A5").HorizontalAlignment = -4152

$ws.Range("B5").Font.Size = 12.1
$ws.Range("B5").Borders.LineStyle = 1
$ws.Range("B5").Borders.Weight = -4138
$ws.Range("B5").Borders.Color = 13421772

$ws.Range("C5").Font.Size = 12.1
$ws.Range("C5").Borders.LineStyle = 1
$ws.Range("C5").Borders.Weight = -4138
$ws.Range("C5").Borders.Color = 13421772

$ws.Range("D5").Font.Size = 10
$ws.Range("D5").Borders.LineStyle = 1
$ws.Range("D5").Borders.Weight = -4138
$ws.Range("D5").Borders.Color = 13421772
$ws.Range("D5").WrapText = $true

$ws.Range("E5").Font.Size = 10
$ws.Range("E5").Borders.LineStyle = 1
$ws.Range("E5").Borders.Weight = -4138
$ws.Range("E5").Borders.Color = 13421772
$ws.Range("E5").WrapText = $true

$ws.Rows(5).RowHeight = 16.5

# ---------------------------------------------------------------------------
# Row 6
# ---------------------------------------------------------------------------
$ws.Range("A6").Font.Size = 12.1
$ws.Range("A6").Borders.LineStyle = 1
$ws.Range("A6").Borders.Weight = -4138
$ws.Range("A6").Borders.Color = 13421772
$ws.Range("A6").HorizontalAlignment = -4152

$ws.Range("B6").Font.Size = 12.1
$ws.Range("B6").Borders.LineStyle = 1
$ws.Range("B6").Borders.Weight = -4138
$ws.Range("B6").Borders.Color = 13421772

$ws.Range("C6").Font.Size = 12.1
$ws.Range("C6").Borders.LineStyle = 1
$ws.Range("C6").Borders.Weight = -4138
$ws.Range("C6").Borders.Color = 13421772

$ws.Range("D6").Font.Size = 12.1
$ws.Range("D6").Borders.LineStyle = 1
$ws.Range("D6").Borders.Weight = -4138
$ws.Range("D6").Borders.Color = 13421772

$ws.Range("E6").Font.Size = 12.1
$ws.Range("E6").Borders.LineStyle = 1
$ws.Range("E6").Borders.Weight = -4138
$ws.Range("E6").Borders.Color = 13421772
$ws.Range("E6").HorizontalAlignment = -4152

$ws.Rows(6).RowHeight = 16.5

# ---------------------------------------------------------------------------
# Row 7
# ---------------------------------------------------------------------------
$ws.Range("A7").Font.Size = 12.1
$ws.Range("A7").Borders.LineStyle = 1
$ws.Range("A7").Borders.Weight = -4138
$ws.Range("A7").Borders.Color = 13421772
$ws.Range("A7").HorizontalAlignment = -4152

$ws.Range("B7").Font.Size = 12.1
$ws.Range("B7").Borders.LineStyle = 1
$ws.Range("B7").Borders.Weight = -4138
$ws.Range("B7").Borders.Color = 13421772

$ws.Range("C7").Font.Size = 12.1
$ws.Range("C7").Borders.LineStyle = 1
$ws.Range("C7").Borders.Weight = -4138
$ws.Range("C7").Borders.Color = 13421772

$ws.Range("D7").Font.Size = 12.1
$ws.Range("D7").Borders.LineStyle = 1
$ws.Range("D7").Borders.Weight = -4138
$ws.Range("D7").Borders.Color = 13421772
$ws.Range("D7").WrapText = $true

$ws.Range("E7").Font.Size = 12.1
$ws.Range("E7").Borders.LineStyle = 1
$ws.Range("E7").Borders.Weight = -4138
$ws.Range("E7").Borders.Color = 13421772
$ws.Range("E7").HorizontalAlignment = -4152

$ws.Rows(7).RowHeight = 48

# ---------------------------------------------------------------------------
# Row 8
# ---------------------------------------------------------------------------
$ws.Range("A8").Font.Size = 12.1
$ws.Range("A8").Font.Color = 0
$ws.Range("A8").Borders.LineStyle = 1
$ws.Range("A8").Borders.Weight = -4138
$ws.Range("A8").Borders.Color = 13421772
$ws.Range("A8").HorizontalAlignment = -4152

$ws.Range("B8").Font.Size = 12.1
$ws.Range("B8").Font.Color = 0
$ws.Range("B8").Borders.LineStyle = 1
$ws.Range("B8").Borders.Weight = -4138
$ws.Range("B8").Borders.Color = 13421772

$ws.Range("C8").Font.Size = 12.1
$ws.Range("C8").Font.Color = 0
$ws.Range("C8").Borders.LineStyle = 1
$ws.Range("C8").Borders.Weight = -4138
$ws.Range("C8").Borders.Color = 13421772

$ws.Range("D8").Font.Size = 12.1
$ws.Range("D8").Font.Color = 0
$ws.Range("D8").Borders.LineStyle = 1
$ws.Range("D8").Borders.Weight = -4138
$ws.Range("D8").Borders.Color = 13421772

$ws.Range("E8").Font.Size = 12.1
$ws.Range("E8").Font.Color = 0
$ws.Range("E8").Borders.LineStyle = 1
$ws.Range("E8").Borders.Weight = -4138
$ws.Range("E8").Borders.Color = 13421772
$ws.Range("E8").HorizontalAlignment = -4152

$ws.Rows(8).RowHeight = 16.5

# ---------------------------------------------------------------------------
# Row 9
# ---------------------------------------------------------------------------
$ws.Range("A9").Font.Size = 12.1
$ws.Range("A9").Font.Color = 0
$ws.Range("A9").Borders.LineStyle = 1
$ws.Range("A9").Borders.Weight = -4138
$ws.Range("A9").Borders.Color = 13421772
$ws.Range("A9").HorizontalAlignment = -4152

$ws.Range("B9").Font.Size = 12.1
$ws.Range("B9").Font.Color = 0
$ws.Range("B9").Borders.LineStyle = 1
$ws.Range("B9").Borders.Weight = -4138
$ws.Range("B9").Borders.Color = 13421772

$ws.Range("C9").Font.Size = 12.1
$ws.Range("C9").Borders.LineStyle = 1
$ws.Range("C9").Borders.Weight = -4138
$ws.Range("C9").Borders.Color = 13421772
$ws.Range("C9").WrapText = $true

$ws.Range("D9").Font.Size = 12.1
$ws.Range("D9").Font.Color = 0
$ws.Range("D9").Borders.LineStyle = 1
$ws.Range("D9").Borders.Weight = -4138
$ws.Range("D9").Borders.Color = 13421772
$ws.Range("D9").WrapText = $true

$ws.Range("E9").Font.Size = 12.1
$ws.Range("E9").Font.Color = 0
$ws.Range("E9").Borders.LineStyle = 1
$ws.Range("E9").Borders.Weight = -4138
$ws.Range("E9").Borders.Color = 13421772
$ws.Range("E9").HorizontalAlignment = -4152

$ws.Rows(9).RowHeight = 48

# ---------------------------------------------------------------------------
# Selection moved to C3 (matches the author's last-edited cell)
# ---------------------------------------------------------------------------
$ws.Range("C3").Select()
